# End of Day 7. Mostly complete
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : Las Pistol ---
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Formula = "=10/B2"
$ws.Range("D2").Formula = "=40/B2"
$ws.Range("E2").Formula = "=8/B2"
$ws.Range("F2").Formula = "=35/B2"
$ws.Range("G2").Formula = "=8/B2"
$ws.Range("H2").Formula = "=25/B2"
$ws.Range("K2").Value = 2

# --- Row 3 : Bolt Pistol ---
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Formula = "=20/B3"
$ws.Range("D3").Formula = "=50/B3"
$ws.Range("E3").Formula = "=15/B3"
$ws.Range("F3").Formula = "=35/B3"
$ws.Range("G3").Formula = "=5/B3"
$ws.Range("H3").Formula = "=25/B3"
$ws.Range("K3").Value = 7

# --- Row 4 : Machine Las Pistol (values unchanged, only K4 added) ---
$ws.Range("K4").Value = 161

# --- Row 5 : Frag Pistol ---
$ws.Range("B5").Value = 0.8
$ws.Range("C5").Formula = "=15/B5"
$ws.Range("D5").Formula = "=60/B5"
$ws.Range("E5").Formula = "=3/B5"
$ws.Range("F5").Formula = "=40/B5"
$ws.Range("G5").Formula = "=3/B5"
$ws.Range("H5").Formula = "=10/B5"
$ws.Range("K5").Value = 180

# --- Row 6 : Frag Pistol (All) ---
$ws.Range("B6").Value = 0.8
$ws.Range("C6").Formula = "=8*15/B6"
$ws.Range("D6").Formula = "=8*60/B6"
$ws.Range("E6").Formula = "=8*8/B6"
$ws.Range("F6").Formula = "=8*40/B6"
$ws.Range("G6").Formula = "=8*3/B6"
$ws.Range("H6").Formula = "=8*10/B6"

# --- Row 7 : Bolt Thrower ---
$ws.Range("B7").Value = 0.25
$ws.Range("C7").Formula = "=40/B7"
$ws.Range("D7").Formula = "=60/B7"
$ws.Range("G7").Formula = "=15/B7"
$ws.Range("H7").Formula = "=50/B7"
$ws.Range("I7").Formula = "=10/B7"
$ws.Range("J7").Formula = "=30/B7"
$ws.Range("K7").Value = 69

# --- Row 8 : Assault Thrower ---
$ws.Range("D8").Formula = "=45/B8"
$ws.Range("G8").Formula = "=10/B8"
$ws.Range("H8").Formula = "=35/B8"
$ws.Range("K8").Value = 152

# --- Row 9 : Plasma Pistol ---
$ws.Range("B9").Value = 0.3
$ws.Range("C9").Formula = "=45/B9"
$ws.Range("D9").Formula = "=70/B9"
$ws.Range("G9").Formula = "=30/B9"
$ws.Range("H9").Formula = "=60/B9"
$ws.Range("I9").Formula = "=10/B9"
$ws.Range("J9").Formula = "=30/B9"
$ws.Range("K9").Value = 80

# --- Row 10 : Rail Gun ---
$ws.Range("C10").Formula = "=125/B10"
$ws.Range("D10").Formula = "=350/B10"
$ws.Range("K10").Value = 22
$ws.Range("L10").Value = 60

# --- Row 11 : Melee (new row) ---
$ws.Range("A11").Value = "Melee"
$ws.Range("B11").Value = 0.4
$ws.Range("C11").Formula = "=100/B11"
$ws.Range("D11").Formula = "=200/B11"
$ws.Range("K11").Value = 193

# --- Selection / active cell ---
$ws.Range("K4").Select() | Out-Null
